$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B10 changes from text "5" to a real numeric value 5
$ws.Cells.Item(10, 2).Value = 5

# New row 11 of annotation data
$ws.Cells.Item(11, 1).Value = "Ying Tang"

# B11 must stay a text value "3" (not numeric), so force text formatting
# before assigning, then reset the cell style back to Normal so no extra
# style index is left behind on the cell.
$ws.Cells.Item(11, 2).NumberFormat = "@"
$ws.Cells.Item(11, 2).Value = "3"
$ws.Cells.Item(11, 2).Style = "Normal"

$ws.Cells.Item(11, 3).Value = "无"
$ws.Cells.Item(11, 4).Value = "QSN"
$ws.Cells.Item(11, 5).Value = "RES"
$ws.Cells.Item(11, 6).Value = "a5228610-fe6d-4383-b598-a7c34c3b8714"
$ws.Cells.Item(11, 7).Value = "HyRnez-RW_annotated.xlsx"
$ws.Cells.Item(11, 8).Value = "Why is this result not compared to in Table 1?"
